# [Kadastro App] Yeni kayit eklendi: 2968
#
# Adds a new record (row 36) to both the "Kayitlar" master sheet and the
# "Erdemli" district sheet, mirroring the structure of the existing rows.
# Columns: Kayit No | Tarih | Birim | Parsel Sayisi | Is | Personeller

$wb = $excel.ActiveWorkbook

function Add-KayitRow($ws, $Row, $KayitNo, $Tarih, $Birim, $ParselSayisi, $Is, $Personeller) {
    # Force text storage for the columns that would otherwise be
    # auto-detected as a number/date by Excel, so the new row keeps the
    # same "number/date stored as text" convention used by every other
    # row already present in this sheet.
    $ws.Range("A$Row").NumberFormat = "@"
    $ws.Range("B$Row").NumberFormat = "@"
    $ws.Range("D$Row").NumberFormat = "@"

    $ws.Range("A$Row").Value = $KayitNo
    $ws.Range("B$Row").Value = $Tarih
    $ws.Range("C$Row").Value = $Birim
    $ws.Range("D$Row").Value = $ParselSayisi
    $ws.Range("E$Row").Value = $Is
    $ws.Range("F$Row").Value = $Personeller
}

$KayitNo = "2968"
$Tarih = "2025-09-10"
$Birim = "Erdemli"
$ParselSayisi = "1"
$Is = "ÇAP"
$Personeller = "AYHAN KARADAYI (K.Teknisyeni)"

# "Kayitlar" is the master log of every record; "Erdemli" mirrors the
# records belonging to that Birim. Both need the new row appended as row 36.
foreach ($sheetName in @("Kayitlar", "Erdemli")) {
    $ws = $wb.Worksheets.Item($sheetName)
    Add-KayitRow $ws 36 $KayitNo $Tarih $Birim $ParselSayisi $Is $Personeller
}
